# feat: add 2022-Q4 data
#
# 1. Insert the new "2022-Q4" quarter row at the top of the "总计" (totals)
#    sheet's data, pushing the existing quarter rows down by one.
# 2. Add a brand-new "2022-Q4" worksheet (positioned right after "总计",
#    i.e. as the 2nd tab) holding the per-fund breakdown for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: "总计" sheet - insert new top row for 2022-Q4, shifting rows down.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 8
$totals.Range("D2").Value = 0.09

# The row-insert shifted the old rows down but kept their old index values in
# column A - renumber them (0-based) to match their new row position.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# Part 2: brand-new "2022-Q4" worksheet with the per-fund breakdown.
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($refSheet)
$q4.Name = "2022-Q4"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$data = @(
  @(0, '006048', '长城中证500指数增强A', '2.20', '94.03', '1.18', '0.0260', 9),
  @(1, '159617', '华夏中证智选500价值稳健策略ETF', '1.74', '97.09', '1.45', '0.0252', 5),
  @(2, '009658', '汇丰晋信中小盘低波动策略股票A', '0.85', '92.42', '1.47', '0.0125', 5),
  @(3, '007413', '长城中证500指数增强C', '0.98', '94.03', '1.18', '0.0116', 9),
  @(4, '159990', '银华巨潮小盘价值ETF', '0.78', '97.02', '1.07', '0.0083', 9),
  @(5, '006346', '安信量化优选股票A', '0.27', '90.65', '1.11', '0.0030', 9),
  @(6, '006347', '安信量化优选股票C', '0.14', '90.65', '1.11', '0.0016', 9),
  @(7, '009775', '汇丰晋信中小盘低波动策略股票C', '0.04', '92.42', '1.47', '0.0006', 5)
)

$r = 2
foreach ($row in $data) {
    $q4.Range("A$r").Value = $row[0]

    # Columns B-G hold numeric-looking fund codes / figures that must stay
    # TEXT (leading zeros, trailing zeros must be preserved) - force the
    # text number format before assigning so the COM layer doesn't coerce
    # them into numbers.
    $q4.Range("B$r").NumberFormat = "@"
    $q4.Range("B$r").Value = $row[1]
    $q4.Range("C$r").Value = $row[2]
    $q4.Range("D$r").NumberFormat = "@"
    $q4.Range("D$r").Value = $row[3]
    $q4.Range("E$r").NumberFormat = "@"
    $q4.Range("E$r").Value = $row[4]
    $q4.Range("F$r").NumberFormat = "@"
    $q4.Range("F$r").Value = $row[5]
    $q4.Range("G$r").NumberFormat = "@"
    $q4.Range("G$r").Value = $row[6]

    $q4.Range("H$r").Value = $row[7]
    $r++
}
